$d = $word.ActiveDocument

# --- Change 1: remove the stray _GoBack bookmark around the "=======" paragraph ---
$d.Bookmarks("_GoBack").Delete()

# --- Change 2: strike-through the "phần đọc thu chi ngan hàng..." paragraph ---
$p = $d.Paragraphs(14)
$p.Range.Font.StrikeThrough = 1

# --- Change 3: wrap the final image paragraph with a fresh _GoBack bookmark ---
$imgPara = $d.Paragraphs(15)
$d.Bookmarks.Add("_GoBack", $imgPara.Range)
